# The document has two logos that were embedded with swapped names:
#   - the Pearson/Edexcel logo (in both footers) was saved as "image1.png"
#     and needs to become "image2.png"
#   - the BTEC logo (in both headers) was saved as "image2.jpg"
#     and needs to become "image1.jpg"
#
# Rename each InlineShape via the standard Word object model
# (InlineShape.Name), reaching every header/footer story of the one
# section in this document.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-LogoShape($headerFooter, [string]$newName) {
    $shape = $headerFooter.Range.InlineShapes.Item(1)
    # Re-seat the shape through Selection before renaming - some stories
    # (notably footers) need a fresh handle or the rename is rejected.
    $null = $shape.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# Headers (BTEC logo): image2.jpg -> image1.jpg
Rename-LogoShape $sec.Headers(1) "image1.jpg"
Rename-LogoShape $sec.Headers(2) "image1.jpg"

# Footers (Pearson/Edexcel logo): image1.png -> image2.png
Rename-LogoShape $sec.Footers(1) "image2.png"
Rename-LogoShape $sec.Footers(2) "image2.png"
